$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 75) with values pulled from Adafruit IO.
# Column C holds a numeric-looking value ("25") that must stay text like
# the rest of the sheet (all cells are stored as strings), so we enter it
# with a leading apostrophe and then reset the style back to the default
# "Normal" style to avoid leaving a stray quote-prefix format behind.
$ws.Range("A75").Value = "2024-09-25T18:06:40Z"
$ws.Range("B75").Value = "temperature"
$ws.Range("C75").Value = "'25"
$ws.Range("C75").Style = "Normal"
$ws.Range("D75").Value = "N/A"
$ws.Range("E75").Value = "N/A"
$ws.Range("F75").Value = "N/A"
